{"js": "// Updated SCPD info about signing up for graders.\n//\n// Replaces the old sentence:\n//   \"or attend a weekly discussion section.  Instead, your grades for other\n//    aspects of the course will be weighted slightly higher to compensate.\"\n// with a longer explanation that tells SCPD students to register at the\n// SCPD signup link, and moves the (Word-managed) \"_GoBack\" bookmark to mark\n// this paragraph as the most recent edit location.\n\nconst OLD_TEXT =\n  \"or attend a weekly discussion section.  Instead, your grades for other \" +\n  \"aspects of the course will be weighted slightly higher to compensate.\";\n\n// A private-use-area character that cannot appear in real document text;\n// used purely as a scratch marker so we can re-locate the exact bookmark\n// insertion point after the full sentence has been inserted. It is removed\n// again before the script finishes.\nconst MARKER = \"\\uE000\";\n\nconst NEW_TEXT =\n  \"or attend a weekly discussion section.  \" +\n  \"Instead, you should register at the SCPD signup link listed in the \\u201cSections\\u201d \" +\n  \"dropdown\" +\n  \" on the course website.  This will allow us to assign you a \" +\n  \"grader\" +\n  \" for \" +\n  MARKER +\n  \"your assignments.  Because you are not attending section, y\" +\n  \"our grades for other aspects of the course will be weighted slightly higher to compensate.\";\n\n// 1. Locate the sentence that needs replacing and swap in the full new text.\nconst body = context.document.body;\nconst found = body.search(OLD_TEXT, { matchCase: true });\nfound.load(\"items\");\nawait context.sync();\n\nif (found.items.length === 0) {\n  throw new Error(\"Could not find target sentence to replace.\");\n}\n\nconst replaced = found.items[0].insertText(NEW_TEXT, Word.InsertLocation.replace);\nawait context.sync();\n\n// 2. Italicize just the word \"grader\" that precedes \"for your assignments\".\n//    Re-search (rather than reusing a sub-range of the just-inserted text)\n//    so the italic formatting doesn't leak into neighboring runs; use\n//    intersectWithOrNullObject against the replaced range so we pick the\n//    right occurrence even though \"grader\" also appears elsewhere in the\n//    document.\nconst graderSearch = context.document.body.search(\"grader\", { matchCase: true });\ngraderSearch.load(\"items\");\nawait context.sync();\n\nlet graderRange = null;\nfor (const candidate of graderSearch.items) {\n  const intersection = candidate.intersectWithOrNullObject(replaced);\n  intersection.load(\"isNullObject\");\n  await context.sync();\n  if (!intersection.isNullObject) {\n    graderRange = candidate;\n    break;\n  }\n}\nif (!graderRange) {\n  throw new Error(\"Could not find the newly-inserted 'grader' run.\");\n}\ngraderRange.font.italic = true;\nawait context.sync();\n\n// 3. Move the \"_GoBack\" bookmark: Word keeps this bookmark at the location\n//    of the user's most recent edit, so it needs to move from its old home\n//    (inside the \"Exam Review\" section) to right after \" for \" / before\n//    \"your assignments\" in the text we just inserted.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst markerSearch = context.document.body.search(MARKER, { matchCase: true });\nmarkerSearch.load(\"items\");\nawait context.sync();\nif (markerSearch.items.length === 0) {\n  throw new Error(\"Could not find the bookmark placeholder marker.\");\n}\nconst markerRange = markerSearch.items[0];\nmarkerRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// Remove the placeholder marker character now that the bookmark sits in its\n// place.\nmarkerRange.insertText(\"\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Updated SCPD info about signing up for graders.\n#\n# Replaces the old sentence:\n#   \"or attend a weekly discussion section.  Instead, your grades for other\n#    aspects of the course will be weighted slightly higher to compensate.\"\n# with a longer explanation that tells SCPD students to register at the\n# SCPD signup link, and moves the (Word-managed) \"_GoBack\" bookmark to mark\n# this paragraph as the most recent edit location.\n\n$d = $word.ActiveDocument\n\n$oldText = \"or attend a weekly discussion section.  Instead, your grades for other aspects of the course will be weighted slightly higher to compensate.\"\n\n$quoteOpen  = [char]0x201C\n$quoteClose = [char]0x201D\n$anchorForText = \" for \"\n\n$newText = \"or attend a weekly discussion section.  \" +\n  \"Instead, you should register at the SCPD signup link listed in the $quoteOpen\" + \"Sections\" + \"$quoteClose \" +\n  \"dropdown\" +\n  \" on the course website.  This will allow us to assign you a \" +\n  \"grader\" +\n  $anchorForText +\n  \"your assignments.  Because you are not attending section, y\" +\n  \"our grades for other aspects of the course will be weighted slightly higher to compensate.\"\n\n# 1. Locate the sentence that needs replacing and swap in the full new text.\n$rng = $d.Content\n$rng.Find.Execute($oldText) | Out-Null\nif (-not $rng.Find.Found) {\n  throw \"Could not find target sentence to replace.\"\n}\n$startPos = $rng.Start\n$rng.Text = $newText\n$endPos = $startPos + $newText.Length\n\n# 2. Italicize just the word \"grader\" that precedes \"for your assignments\".\n#    Scope the Find to the range we just inserted (by absolute character\n#    position) so we don't accidentally match one of the other \"grader\"\n#    occurrences elsewhere in the document.\n$scopedRange = $d.Range($startPos, $endPos)\n$scopedRange.Find.Execute(\"grader\") | Out-Null\nif (-not $scopedRange.Find.Found) {\n  throw \"Could not find the newly-inserted 'grader' run.\"\n}\n$scopedRange.Font.Italic = 1\n\n# 3. Move the \"_GoBack\" bookmark: Word keeps this bookmark at the location\n#    of the user's most recent edit, so it needs to move from its old home\n#    (inside the \"Exam Review\" section) to right after \" for \" / before\n#    \"your assignments\" in the text we just inserted.\nif ($d.Bookmarks.Exists('_GoBack')) {\n  $d.Bookmarks.Item('_GoBack').Delete()\n}\n\n$bmScope = $d.Range($startPos, $endPos)\n$bmScope.Find.Execute(\"$anchorForText\" + \"your assignments\") | Out-Null\nif (-not $bmScope.Find.Found) {\n  throw \"Could not find insertion point for the _GoBack bookmark.\"\n}\n$bmPos = $bmScope.Start + $anchorForText.Length\n$bmRange = $d.Range($bmPos, $bmPos)\n$d.Bookmarks.Add('_GoBack', $bmRange) | Out-Null\n"}
